# Auto-generated: full refresh of country COVID stats table + re-sort by Casos totales desc
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "Datos actualizados a 9 de Abril de 2020 a las 16:52"

$data = @(
    ,@("Estados Unidos", 435803, 876, 22941, 397997, 9281, 77, 14865)
    ,@("España", 152446, 4226, 52165, 85043, 7371, 446, 15238)
    ,@("Italia", 139422, 0, 26491, 95262, 3693, 0, 17669)
    ,@("Alemania", 113615, 319, 46300, 64966, 4895, 0, 2349)
    ,@("Francia", 112950, 0, 21254, 80827, 7148, 0, 10869)
    ,@("China", 81865, 63, 77370, 1160, 176, 2, 3335)
    ,@("Iran", 66220, 1634, 32309, 29801, 3918, 117, 4110)
    ,@("Reino Unido", 60733, 0, 135, 53501, 1559, 0, 7097)
    ,@("Turquia", 38226, 0, 1846, 35568, 1492, 0, 812)
    ,@("Belgica", 24983, 1580, 5164, 17296, 1285, 283, 2523)
    ,@("Suiza", 23799, 519, 9900, 12962, 386, 42, 937)
    ,@("Paises Bajos", 21762, 1213, 250, 19116, 1424, 148, 2396)
    ,@("Canada", 19438, 0, 4653, 14350, 518, 8, 435)
    ,@("Brasil", 16275, 87, 173, 15276, 296, 6, 826)
    ,@("Portugal", 13956, 815, 205, 13342, 241, 29, 409)
    ,@("Austria", 13163, 221, 5240, 7628, 266, 22, 295)
    ,@("Corea del Sur", 10423, 39, 6973, 3246, 55, 4, 204)
    ,@("Rusia", 10131, 1459, 698, 9357, 8, 13, 76)
    ,@("Israel", 9755, 351, 864, 8812, 165, 6, 79)
    ,@("Suecia", 9141, 722, 205, 8143, 719, 106, 793)
    ,@("India", 6237, 321, 569, 5482, 0, 8, 186)
    ,@("Noruega", 6160, 118, 32, 6024, 78, 3, 104)
    ,@("Australia", 6104, 52, 2987, 3066, 81, 1, 51)
    ,@("Irlanda", 6074, 0, 25, 5814, 165, 0, 235)
    ,@("Chile", 5972, 426, 1274, 4641, 360, 9, 57)
    ,@("Dinamarca", 5635, 233, 1736, 3662, 120, 19, 237)
    ,@("Polonia", 5341, 136, 284, 4893, 160, 5, 164)
    ,@("Chequia", 5335, 23, 243, 4988, 96, 5, 104)
    ,@("Rumania", 5202, 441, 647, 4326, 178, 9, 229)
    ,@("Japon", 4667, 0, 632, 3941, 99, 0, 94)
    ,@("Pakistan", 4462, 199, 572, 3827, 31, 2, 63)
    ,@("Ecuador", 4450, 0, 140, 4068, 146, 0, 242)
    ,@("Peru", 4342, 0, 1333, 2888, 113, 0, 121)
    ,@("Malasia", 4228, 109, 1608, 2553, 72, 2, 67)
    ,@("Filipinas", 4076, 206, 124, 3749, 1, 21, 203)
    ,@("Indonesia", 3293, 337, 252, 2761, 0, 40, 280)
    ,@("Arabia Saudita", 3287, 355, 666, 2577, 41, 3, 44)
    ,@("Mexico", 3181, 396, 633, 2374, 89, 33, 174)
    ,@("Luxemburgo", 3034, 0, 500, 2488, 34, 0, 46)
    ,@("Serbia", 2867, 201, 118, 2683, 127, 1, 66)
    ,@("Emiratos Arabes Unidos", 2659, 0, 239, 2408, 1, 0, 12)
    ,@("Finlandia", 2605, 118, 300, 2263, 82, 2, 42)
    ,@("Panama", 2528, 0, 16, 2449, 101, 0, 63)
    ,@("Tailandia", 2423, 54, 940, 1451, 61, 2, 32)
    ,@("Catar", 2376, 166, 206, 2164, 37, 0, 6)
    ,@("Republica Dominicana", 2111, 0, 50, 1953, 147, 0, 108)
    ,@("Colombia", 2054, 0, 123, 1876, 76, 0, 55)
    ,@("Singapur", 1910, 287, 406, 1498, 29, 0, 6)
    ,@("Ucrania", 1892, 224, 45, 1790, 33, 5, 57)
    ,@("Grecia", 1884, 0, 269, 1532, 84, 0, 83)
    ,@("Sudafrica", 1845, 0, 95, 1732, 7, 0, 18)
    ,@("Argentina", 1795, 0, 358, 1367, 96, 5, 70)
    ,@("Islandia", 1616, 0, 633, 977, 13, 0, 6)
    ,@("Argelia", 1572, 0, 237, 1130, 46, 0, 205)
    ,@("Egipto", 1560, 0, 305, 1152, 0, 0, 103)
    ,@("Bielorrusia", 1486, 420, 139, 1331, 55, 3, 16)
    ,@("Croacia", 1407, 64, 219, 1168, 34, 1, 20)
    ,@("Marruecos", 1346, 71, 103, 1147, 1, 3, 96)
    ,@("Moldavia", 1289, 115, 50, 1210, 80, 2, 29)
    ,@("Nueva Zelanda", 1239, 29, 317, 921, 4, 0, 1)
    ,@("Estonia", 1207, 22, 83, 1100, 9, 0, 24)
    ,@("Irak", 1202, 0, 452, 681, 0, 0, 69)
    ,@("Eslovenia", 1124, 33, 128, 953, 34, 3, 43)
    ,@("Hungria", 980, 85, 96, 818, 17, 8, 66)
    ,@("Hong Kong", 974, 13, 293, 677, 14, 0, 4)
    ,@("Lituania", 955, 43, 8, 931, 21, 1, 16)
    ,@("Azerbaiyan", 926, 104, 101, 816, 27, 1, 9)
    ,@("Armenia", 921, 40, 138, 773, 30, 1, 10)
    ,@("Kuwait", 910, 55, 111, 798, 22, 0, 1)
    ,@("Barein", 855, 32, 495, 355, 3, 0, 5)
    ,@("Bosnia y Herzegovina", 841, 37, 95, 711, 4, 1, 35)
    ,@("Kazajistan", 764, 37, 60, 697, 21, 0, 7)
    ,@("Camerun", 730, 0, 60, 660, 0, 0, 10)
    ,@("Crucero", 712, 0, 619, 82, 10, 0, 11)
    ,@("Eslovaquia", 701, 19, 23, 676, 5, 0, 2)
    ,@("Republica de Macedonia", 663, 46, 37, 596, 15, 1, 30)
    ,@("Tunez", 628, 0, 25, 579, 78, 0, 24)
    ,@("Bulgaria", 618, 25, 48, 546, 27, 0, 24)
    ,@("Letonia", 589, 12, 16, 571, 3, 0, 2)
    ,@("Principado de Andorra", 583, 19, 58, 500, 17, 2, 25)
    ,@("Uzbekistan", 582, 37, 38, 541, 8, 0, 3)
    ,@("Libano", 582, 6, 67, 496, 29, 0, 19)
    ,@("Republica de Chipre", 526, 0, 52, 465, 13, 0, 9)
    ,@("Costa Rica", 502, 0, 29, 470, 15, 0, 3)
    ,@("Afganistan", 484, 40, 32, 437, 0, 1, 15)
    ,@("Cuba", 457, 0, 27, 418, 15, 0, 12)
    ,@("Oman", 457, 38, 109, 346, 3, 0, 2)
    ,@("Uruguay", 456, 0, 192, 257, 14, 0, 7)
    ,@("Burkina Faso", 414, 0, 134, 257, 0, 0, 23)
    ,@("Albania", 409, 9, 165, 222, 7, 0, 22)
    ,@("Costa de Marfil", 384, 0, 48, 333, 0, 0, 3)
    ,@("Taiwan", 380, 1, 80, 295, 0, 0, 5)
    ,@("Reunion", 362, 0, 40, 322, 4, 0, 0)
    ,@("Jordania", 358, 0, 150, 202, 5, 0, 6)
    ,@("Honduras", 343, 31, 6, 314, 10, 1, 23)
    ,@("Niger", 342, 0, 28, 303, 0, 0, 11)
    ,@("Malta", 337, 38, 16, 319, 4, 1, 2)
    ,@("Banglades", 330, 112, 33, 276, 1, 1, 21)
    ,@("Mauricio", 314, 41, 23, 284, 3, 0, 7)
    ,@("Ghana", 313, 0, 34, 273, 2, 0, 6)
    ,@("San Marino", 308, 0, 45, 229, 14, 0, 34)
    ,@("Kirguistan", 280, 10, 35, 241, 5, 0, 4)
    ,@("Nigeria", 276, 0, 44, 226, 2, 0, 6)
    ,@("Bolivia", 264, 54, 2, 244, 3, 3, 18)
    ,@("Estado de Palestina", 263, 0, 44, 218, 0, 0, 1)
    ,@("Vietnam", 255, 4, 128, 127, 8, 0, 0)
    ,@("Montenegro", 252, 4, 4, 246, 7, 0, 2)
    ,@("Senegal", 250, 6, 123, 125, 1, 0, 2)
    ,@("Georgia", 214, 3, 50, 161, 6, 0, 3)
    ,@("Sri Lanka", 190, 1, 49, 134, 5, 0, 7)
    ,@("Kenia", 184, 5, 12, 165, 2, 1, 7)
    ,@("Mayotte", 184, 0, 22, 160, 3, 0, 2)
    ,@("Islas Feroe", 184, 0, 136, 48, 0, 0, 0)
    ,@("Consejo Danes para los Refugiados", 180, 0, 9, 153, 0, 0, 18)
    ,@("Isla de Man", 171, 13, 88, 82, 10, 0, 1)
    ,@("Venezuela", 167, 0, 65, 93, 6, 0, 9)
    ,@("Guinea", 164, 0, 5, 159, 0, 0, 0)
    ,@("Martinica", 154, 0, 50, 98, 19, 0, 6)
    ,@("Guadalupe", 141, 0, 43, 90, 13, 0, 8)
    ,@("Republica de Yibuti", 135, 0, 25, 110, 0, 0, 0)
    ,@("Brunei", 135, 0, 92, 42, 3, 0, 1)
    ,@("Paraguay", 124, 5, 18, 101, 1, 0, 5)
    ,@("Gibraltar", 120, 0, 60, 60, 0, 0, 0)
    ,@("Camboya", 119, 2, 63, 56, 1, 0, 0)
    ,@("Ruanda", 110, 0, 7, 103, 0, 0, 0)
    ,@("Trinidad yTobago", 109, 2, 1, 100, 0, 0, 8)
    ,@("El Salvador", 103, 10, 10, 88, 3, 0, 5)
    ,@("Guatemala", 95, 8, 17, 75, 3, 0, 3)
    ,@("Madagascar", 93, 0, 11, 82, 1, 0, 0)
    ,@("Guayana Francesa", 83, 0, 43, 40, 1, 0, 0)
    ,@("Monaco", 81, 0, 4, 76, 4, 0, 1)
    ,@("Liechtenstein", 78, 0, 55, 22, 0, 0, 1)
    ,@("Aruba", 77, 0, 14, 63, 0, 0, 0)
    ,@("Mali", 74, 15, 22, 45, 0, 0, 7)
    ,@("Togo", 73, 3, 24, 46, 0, 0, 3)
    ,@("Barbados", 63, 0, 8, 52, 4, 0, 3)
    ,@("Jamaica", 63, 0, 10, 49, 0, 0, 4)
    ,@("Congo", 60, 15, 5, 50, 0, 0, 5)
    ,@("Etiopia", 56, 1, 4, 50, 2, 0, 2)
    ,@("Uganda", 53, 0, 0, 53, 0, 0, 0)
    ,@("Polinesia Francesa", 51, 0, 0, 51, 0, 0, 0)
    ,@("Islas Caimanes", 45, 0, 6, 38, 0, 0, 1)
    ,@("Macao", 45, 0, 10, 35, 1, 0, 0)
    ,@("San Martin (Parte Holandesa)", 43, 3, 1, 36, 2, 0, 6)
    ,@("Bahamas", 40, 0, 5, 28, 1, 0, 7)
    ,@("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
    ,@("Zambia", 39, 0, 7, 31, 1, 0, 1)
    ,@("Bermudas", 39, 0, 23, 13, 0, 0, 3)
    ,@("Guyana", 37, 0, 8, 23, 4, 0, 6)
    ,@("Gabon", 34, 0, 1, 32, 0, 0, 1)
    ,@("Eritrea", 33, 0, 0, 33, 0, 0, 0)
    ,@("Guinea-Bisau", 33, 0, 0, 33, 0, 0, 0)
    ,@("Guam", 32, 0, 0, 31, 0, 0, 1)
    ,@("San Martin (Parte Francesa)", 32, 0, 7, 23, 6, 0, 2)
    ,@("Liberia", 31, 0, 3, 24, 0, 0, 4)
    ,@("Haiti", 30, 3, 0, 28, 0, 1, 2)
    ,@("Benin", 26, 0, 5, 20, 0, 0, 1)
    ,@("Tanzania", 25, 0, 5, 19, 0, 0, 1)
    ,@("Birmania", 23, 1, 1, 19, 0, 0, 3)
    ,@("Libia", 21, 0, 8, 12, 0, 0, 1)
    ,@("Antigua y Barbuda", 19, 0, 0, 17, 1, 0, 2)
    ,@("Angola", 19, 0, 2, 15, 0, 0, 2)
    ,@("Siria", 19, 0, 4, 13, 0, 0, 2)
    ,@("Maldivas", 19, 0, 13, 6, 0, 0, 0)
    ,@("Nueva Caledonia", 18, 0, 1, 17, 0, 0, 0)
    ,@("Guinea Ecuatorial", 18, 0, 3, 15, 0, 0, 0)
    ,@("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
    ,@("Mozambique", 17, 0, 1, 16, 0, 0, 0)
    ,@("Laos", 16, 1, 0, 16, 0, 0, 0)
    ,@("Namibia", 16, 0, 3, 13, 0, 0, 0)
    ,@("Mongolia", 16, 0, 4, 12, 0, 0, 0)
    ,@("Fiyi", 15, 0, 0, 15, 0, 0, 0)
    ,@("Dominica", 15, 0, 1, 14, 0, 0, 0)
    ,@("Santa Lucia", 14, 0, 1, 13, 0, 0, 0)
    ,@("Sudan", 14, 0, 2, 10, 0, 0, 2)
    ,@("Curazao", 14, 0, 7, 6, 0, 0, 1)
    ,@("Botsuana", 13, 7, 0, 12, 0, 0, 1)
    ,@("Granada", 12, 0, 0, 12, 2, 0, 0)
    ,@("Somalia", 12, 0, 1, 10, 0, 0, 1)
    ,@("Suazilandia", 12, 0, 7, 5, 0, 0, 0)
    ,@("Seychelles", 11, 0, 0, 11, 0, 0, 0)
    ,@("San Cristobal y Nieves", 11, 0, 0, 11, 0, 0, 0)
    ,@("Zimbabue", 11, 0, 0, 8, 0, 0, 3)
    ,@("Groenlandia", 11, 0, 11, 0, 0, 0, 0)
    ,@("Republica del Chad", 10, 0, 2, 8, 0, 0, 0)
    ,@("Surinam", 10, 0, 4, 5, 0, 0, 1)
    ,@("Nepal", 9, 0, 1, 8, 0, 0, 0)
    ,@("Montserrat", 9, 0, 0, 7, 0, 0, 2)
    ,@("Republica de Africa Central", 8, 0, 0, 8, 0, 0, 0)
    ,@("San Vicente y las Granadinas", 8, 0, 1, 7, 0, 0, 0)
    ,@("Malaui", 8, 0, 0, 7, 1, 0, 1)
    ,@("Islas Turcas y Caicos", 8, 0, 0, 7, 0, 0, 1)
    ,@("Belice", 8, 0, 0, 7, 1, 0, 1)
    ,@("Santa Sede", 8, 0, 2, 6, 0, 0, 0)
    ,@("Sierra Leona", 7, 0, 0, 7, 0, 0, 0)
    ,@("Cabo Verde", 7, 0, 1, 5, 0, 0, 1)
    ,@("Mauritania", 7, 1, 2, 4, 0, 0, 1)
    ,@("San Bartolome", 6, 0, 1, 5, 0, 0, 0)
    ,@("Nicaragua", 6, 0, 0, 5, 0, 0, 1)
    ,@("Islas Malvinas", 5, 0, 1, 4, 0, 0, 0)
    ,@("Butan", 5, 0, 2, 3, 0, 0, 0)
    ,@("Sahara Occidental", 4, 0, 0, 4, 0, 0, 0)
    ,@("Santo Tome y Principe", 4, 0, 0, 4, 0, 0, 0)
    ,@("Gambia", 4, 0, 2, 1, 0, 0, 1)
    ,@("Sudan del Sur", 3, 1, 0, 3, 0, 0, 0)
    ,@("Burundi", 3, 0, 0, 3, 0, 0, 0)
    ,@("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0)
    ,@("Anguila", 3, 0, 0, 3, 0, 0, 0)
    ,@("Bonaire, San Eustaquio y Saba", 2, 0, 0, 2, 0, 0, 0)
    ,@("Papua Nueva Guinea", 2, 0, 0, 2, 0, 0, 0)
    ,@("Timor Oriental", 1, 0, 0, 1, 0, 0, 0)
    ,@("San Pedro y Miquelon", 1, 0, 0, 1, 0, 0, 0)
)

$r = 4
foreach ($item in $data) {
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
    $ws.Cells.Item($r, 6).Value = $item[5]
    $ws.Cells.Item($r, 7).Value = $item[6]
    $ws.Cells.Item($r, 8).Value = $item[7]
    $r++
}

Write-Output "Wrote $($data.Length) rows"
